$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 28.479168
$ws.Range("H2").Value = 85.437504
$ws.Range("I2").Value = 0.4446244458164738
$ws.Range("J2").Value = 0.4446244458164738
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 35.31114333333333
$ws.Range("N2").Value = 105.93343
$ws.Range("O2").Value = 0.6187867769880316
$ws.Range("P2").Value = 0.6187867769880316
$ws.Range("Q2").Value = 1005.63198326208
$ws.Range("R2").Value = 9050.687849358721
$ws.Range("S2").Value = 0.2751277277968656
$ws.Range("T2").Value = 0.2751277277968656

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 28.479168
$ws.Range("H3").Value = 85.437504
$ws.Range("I3").Value = 0.4446244458164738
$ws.Range("J3").Value = 0.4446244458164738
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 6.551362
$ws.Range("N3").Value = 19.654086
$ws.Range("O3").Value = 0.1148050103785518
$ws.Range("P3").Value = 0.1148050103785518
$ws.Range("Q3").Value = 186.577339026816
$ws.Range("R3").Value = 1679.196051241344
$ws.Range("S3").Value = 0.05104511411651813
$ws.Range("T3").Value = 0.05104511411651813

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 28.479168
$ws.Range("H4").Value = 85.437504
$ws.Range("I4").Value = 0.4446244458164738
$ws.Range("J4").Value = 0.4446244458164738
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 10.951915
$ws.Range("N4").Value = 32.855745
$ws.Range("O4").Value = 0.191919590955288
$ws.Range("P4").Value = 0.191919590955288
$ws.Range("Q4").Value = 311.90142720672
$ws.Range("R4").Value = 2807.11284486048
$ws.Range("S4").Value = 0.08533214176981926
$ws.Range("T4").Value = 0.08533214176981926

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 28.479168
$ws.Range("H5").Value = 85.437504
$ws.Range("I5").Value = 0.4446244458164738
$ws.Range("J5").Value = 0.4446244458164738
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 4.250702333333333
$ws.Range("N5").Value = 12.752107
$ws.Range("O5").Value = 0.07448862167812857
$ws.Range("P5").Value = 0.07448862167812857
$ws.Range("Q5").Value = 121.056465868992
$ws.Range("R5").Value = 1089.508192820928
$ws.Range("S5").Value = 0.0331194621332709
$ws.Range("T5").Value = 0.0331194621332709

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 18.12667766666667
$ws.Range("H6").Value = 54.380033
$ws.Range("I6").Value = 0.2829985767855128
$ws.Range("J6").Value = 0.2829985767855128
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 35.31114333333333
$ws.Range("N6").Value = 105.93343
$ws.Range("O6").Value = 0.6187867769880316
$ws.Range("P6").Value = 0.6187867769880316
$ws.Range("Q6").Value = 640.073713244799
$ws.Range("R6").Value = 5760.663419203191
$ws.Range("S6").Value = 0.1751157772213075
$ws.Range("T6").Value = 0.1751157772213075

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 18.12667766666667
$ws.Range("H7").Value = 54.380033
$ws.Range("I7").Value = 0.2829985767855128
$ws.Range("J7").Value = 0.2829985767855128
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 6.551362
$ws.Range("N7").Value = 19.654086
$ws.Range("O7").Value = 0.1148050103785518
$ws.Range("P7").Value = 0.1148050103785518
$ws.Range("Q7").Value = 118.7544272516487
$ws.Range("R7").Value = 1068.789845264838
$ws.Range("S7").Value = 0.0324896545449762
$ws.Range("T7").Value = 0.0324896545449762

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 18.12667766666667
$ws.Range("H8").Value = 54.380033
$ws.Range("I8").Value = 0.2829985767855128
$ws.Range("J8").Value = 0.2829985767855128
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 10.951915
$ws.Range("N8").Value = 32.855745
$ws.Range("O8").Value = 0.191919590955288
$ws.Range("P8").Value = 0.191919590955288
$ws.Range("Q8").Value = 198.5218330377317
$ws.Range("R8").Value = 1786.696497339585
$ws.Range("S8").Value = 0.05431297109760428
$ws.Range("T8").Value = 0.05431297109760428

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 18.12667766666667
$ws.Range("H9").Value = 54.380033
$ws.Range("I9").Value = 0.2829985767855128
$ws.Range("J9").Value = 0.2829985767855128
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 4.250702333333333
$ws.Range("N9").Value = 12.752107
$ws.Range("O9").Value = 0.07448862167812857
$ws.Range("P9").Value = 0.07448862167812857
$ws.Range("Q9").Value = 77.05111105328123
$ws.Range("R9").Value = 693.459999479531
$ws.Range("S9").Value = 0.02108017392162488
$ws.Range("T9").Value = 0.02108017392162488

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 11.513346
$ws.Range("H10").Value = 34.540038
$ws.Range("I10").Value = 0.179749460544048
$ws.Range("J10").Value = 0.179749460544048
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 35.31114333333333
$ws.Range("N10").Value = 105.93343
$ws.Range("O10").Value = 0.6187867769880316
$ws.Range("P10").Value = 0.6187867769880316
$ws.Range("Q10").Value = 406.54941085226
$ws.Range("R10").Value = 3658.94469767034
$ws.Range("S10").Value = 0.1112265893553888
$ws.Range("T10").Value = 0.1112265893553888

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 11.513346
$ws.Range("H11").Value = 34.540038
$ws.Range("I11").Value = 0.179749460544048
$ws.Range("J11").Value = 0.179749460544048
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 6.551362
$ws.Range("N11").Value = 19.654086
$ws.Range("O11").Value = 0.1148050103785518
$ws.Range("P11").Value = 0.1148050103785518
$ws.Range("Q11").Value = 75.42809747725201
$ws.Range("R11").Value = 678.8528772952681
$ws.Range("S11").Value = 0.02063613868329853
$ws.Range("T11").Value = 0.02063613868329853

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 11.513346
$ws.Range("H12").Value = 34.540038
$ws.Range("I12").Value = 0.179749460544048
$ws.Range("J12").Value = 0.179749460544048
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 10.951915
$ws.Range("N12").Value = 32.855745
$ws.Range("O12").Value = 0.191919590955288
$ws.Range("P12").Value = 0.191919590955288
$ws.Range("Q12").Value = 126.09318675759
$ws.Range("R12").Value = 1134.83868081831
$ws.Range("S12").Value = 0.03449744294204737
$ws.Range("T12").Value = 0.03449744294204737

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 11.513346
$ws.Range("H13").Value = 34.540038
$ws.Range("I13").Value = 0.179749460544048
$ws.Range("J13").Value = 0.179749460544048
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 4.250702333333333
$ws.Range("N13").Value = 12.752107
$ws.Range("O13").Value = 0.07448862167812857
$ws.Range("P13").Value = 0.07448862167812857
$ws.Range("Q13").Value = 48.939806706674
$ws.Range("R13").Value = 440.458260360066
$ws.Range("S13").Value = 0.01338928956331329
$ws.Range("T13").Value = 0.01338928956331329

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 5.932994999999999
$ws.Range("H14").Value = 17.798985
$ws.Range("I14").Value = 0.09262751685396531
$ws.Range("J14").Value = 0.09262751685396531
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 35.31114333333333
$ws.Range("N14").Value = 105.93343
$ws.Range("O14").Value = 0.6187867769880316
$ws.Range("P14").Value = 0.6187867769880316
$ws.Range("Q14").Value = 209.50083684095
$ws.Range("R14").Value = 1885.50753156855
$ws.Range("S14").Value = 0.05731668261446977
$ws.Range("T14").Value = 0.05731668261446977

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 5.932994999999999
$ws.Range("H15").Value = 17.798985
$ws.Range("I15").Value = 0.09262751685396531
$ws.Range("J15").Value = 0.09262751685396531
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 6.551362
$ws.Range("N15").Value = 19.654086
$ws.Range("O15").Value = 0.1148050103785518
$ws.Range("P15").Value = 0.1148050103785518
$ws.Range("Q15").Value = 38.86919798918999
$ws.Range("R15").Value = 349.8227819027099
$ws.Range("S15").Value = 0.01063410303375897
$ws.Range("T15").Value = 0.01063410303375897

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 5.932994999999999
$ws.Range("H16").Value = 17.798985
$ws.Range("I16").Value = 0.09262751685396531
$ws.Range("J16").Value = 0.09262751685396531
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 10.951915
$ws.Range("N16").Value = 32.855745
$ws.Range("O16").Value = 0.191919590955288
$ws.Range("P16").Value = 0.191919590955288
$ws.Range("Q16").Value = 64.97765693542499
$ws.Range("R16").Value = 584.7989124188249
$ws.Range("S16").Value = 0.01777703514581707
$ws.Range("T16").Value = 0.01777703514581707

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 5.932994999999999
$ws.Range("H17").Value = 17.798985
$ws.Range("I17").Value = 0.09262751685396531
$ws.Range("J17").Value = 0.09262751685396531
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 4.250702333333333
$ws.Range("N17").Value = 12.752107
$ws.Range("O17").Value = 0.07448862167812857
$ws.Range("P17").Value = 0.07448862167812857
$ws.Range("Q17").Value = 25.21939569015499
$ws.Range("R17").Value = 226.974561211395
$ws.Range("S17").Value = 0.0068996960599195
$ws.Range("T17").Value = 0.0068996960599195
